$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the "harvester" column (B) value for every data row with the
# actual harvester, and fill in the new "experimentDesign" column (D)
# with the induction timing used for the experiment.
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 2).Value = "S.GISH"
    $ws.Cells.Item($r, 4).Value = "90minuteInduction"
}

# Mirror the author's final selection (D2:D19) left after filling the
# new column.
[void]$ws.Range("D2:D19").Select()
